$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (column I) and IF (column J), matching the style of
# the existing header row (style index 1, same as B1:H1). Use copy/paste
# so the existing cell style gets reused rather than cloned.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill in the data rows: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 30; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
